$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value2 = 169
$ws.Range("F6").Value2 = 2753
$ws.Range("F8").Value2 = 1636
$ws.Range("F9").Value2 = 7451
$ws.Range("F11").Value2 = 7646
$ws.Range("F13").Value2 = 35
$ws.Range("F14").Value2 = 7
$ws.Range("F15").Value2 = 6162
$ws.Range("F16").Value2 = 3258
$ws.Range("F18").Value2 = 17
$ws.Range("F19").Value2 = 11
$ws.Range("F20").Value2 = 17
$ws.Range("F21").Value2 = 29
$ws.Range("F22").Value2 = 444
$ws.Range("F25").Value2 = 286
$ws.Range("F26").Value2 = 3624
$ws.Range("F31").Value2 = 1091
$ws.Range("F32").Value2 = 63
$ws.Range("F33").Value2 = 20
$ws.Range("F34").Value2 = 2616
$ws.Range("F35").Value2 = 1462
$ws.Range("F36").Value2 = 12
$ws.Range("F37").Value2 = 19
$ws.Range("F39").Value2 = 3262
$ws.Range("F40").Value2 = 164
$ws.Range("F41").Value2 = 243
$ws.Range("F45").Value2 = 1278
$ws.Range("F46").Value2 = 225
$ws.Range("F47").Value2 = 524
$ws.Range("F48").Value2 = 591

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F13").Value2 = 15

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value2 = 120

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value2 = 169
$ws.Range("F8").Value2 = 120
$ws.Range("F9").Value2 = 1636
$ws.Range("F12").Value2 = 7451
$ws.Range("F13").Value2 = 7646
$ws.Range("F15").Value2 = 6162
$ws.Range("F16").Value2 = 3258
$ws.Range("F18").Value2 = 17
$ws.Range("F19").Value2 = 11
$ws.Range("F20").Value2 = 17
$ws.Range("F21").Value2 = 29
$ws.Range("F22").Value2 = 444
$ws.Range("F26").Value2 = 286
$ws.Range("F27").Value2 = 3624
$ws.Range("F33").Value2 = 63
$ws.Range("F34").Value2 = 20
$ws.Range("F35").Value2 = 2616
$ws.Range("F36").Value2 = 1462
$ws.Range("F37").Value2 = 12
$ws.Range("F38").Value2 = 19
$ws.Range("F40").Value2 = 3262
$ws.Range("F41").Value2 = 164
$ws.Range("F42").Value2 = 243
$ws.Range("F47").Value2 = 1278
$ws.Range("F48").Value2 = 225
$ws.Range("F49").Value2 = 524
